# Update the TPM-derived values in the LR-pairs sheet (Gpc3-Cd81).
# The underlying ligand/receptor expression inputs changed (new TPM values),
# which cascades into the derived specificity / edge-weight columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value
$updates = @{
    "G2" = 0.301437
    "H2" = 0.9043110000000001
    "I2" = 0.003123224212368042
    "J2" = 0.003123224212368043
    "M2" = 87.77304733333334
    "N2" = 263.319142
    "O2" = 0.200063793449224
    "P2" = 0.200063793449224
    "Q2" = 26.458044069018
    "R2" = 238.122396621162
    "S2" = 0.0006248440837188155
    "T2" = 0.0006248440837188156

    "G3" = 0.301437
    "H3" = 0.9043110000000001
    "I3" = 0.003123224212368042
    "J3" = 0.003123224212368043
    "O3" = 0.4942765199240737
    "P3" = 0.4942765199240737
    "Q3" = 65.36709976836001
    "R3" = 588.3038979152401
    "S3" = 0.001543736394631882
    "T3" = 0.001543736394631882

    "G4" = 0.301437
    "H4" = 0.9043110000000001
    "I4" = 0.003123224212368042
    "J4" = 0.003123224212368043
    "M4" = 134.100637
    "N4" = 402.301911
    "O4" = 0.3056596866267023
    "P4" = 0.3056596866267022
    "Q4" = 40.422893715369
    "R4" = 363.8060434383211
    "S4" = 0.0009546437340173449
    "T4" = 0.0009546437340173449

    "I5" = 0.7782793322359159
    "J5" = 0.7782793322359159
    "M5" = 87.77304733333334
    "N5" = 263.319142
    "O5" = 0.200063793449224
    "P5" = 0.200063793449224
    "Q5" = 6593.106184551192
    "R5" = 59337.95566096073
    "S5" = 0.1557055155702463
    "T5" = 0.1557055155702463

    "I6" = 0.7782793322359159
    "J6" = 0.7782793322359159
    "O6" = 0.4942765199240737
    "P6" = 0.4942765199240737
    "S6" = 0.3846851998664005
    "T6" = 0.3846851998664005

    "I7" = 0.7782793322359159
    "J7" = 0.7782793322359159
    "M7" = 134.100637
    "N7" = 402.301911
    "O7" = 0.3056596866267023
    "P7" = 0.3056596866267022
    "Q7" = 10073.02088759982
    "R7" = 90657.18798839839
    "S7" = 0.2378886167992692
    "T7" = 0.2378886167992691

    "G8" = 21.09786333333333
    "H8" = 63.29359
    "I8" = 0.2185974435517159
    "J8" = 0.218597443551716
    "M8" = 87.77304733333334
    "N8" = 263.319142
    "O8" = 0.200063793449224
    "P8" = 0.200063793449224
    "Q8" = 1851.823756988864
    "R8" = 16666.41381289978
    "S8" = 0.0437334337952589
    "T8" = 0.04373343379525891

    "G9" = 21.09786333333333
    "H9" = 63.29359
    "I9" = 0.2185974435517159
    "J9" = 0.218597443551716
    "O9" = 0.4942765199240737
    "P9" = 0.4942765199240737
    "Q9" = 4575.105701719512
    "R9" = 41175.9513154756
    "S9" = 0.1080475836630413
    "T9" = 0.1080475836630413

    "G10" = 21.09786333333333
    "H10" = 63.29359
    "I10" = 0.2185974435517159
    "J10" = 0.218597443551716
    "M10" = 134.100637
    "N10" = 402.301911
    "O10" = 0.3056596866267023
    "P10" = 0.3056596866267022
    "Q10" = 2829.236912338944
    "R10" = 25463.13221105049
    "S10" = 0.06681642609341573
    "T10" = 0.06681642609341573
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
